# "Fixing Case sensitive issues"
# 1) Two existing food_item names had inconsistent capitalisation -> lower-cased.
# 2) A batch of new rows (303-310) of food items were appended to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Case fixes on existing rows ---------------------------------------
$ws.Cells.Item(78, 1).Value  = "kellogg's special k red berries cereal"
$ws.Cells.Item(112, 1).Value = "kirkland salted mixed nuts"

# --- 2. Append new rows 303-310 --------------------------------------------
$newRows = @(
    @(303, "planet oat creamer",                        "1 tbsp",    25,                 0,                    1,                    4),
    @(304, "veer ginger garlic paste",                   "1 tbsp",    170,                4,                    8,                    12),
    @(305, "mtr tomato rice powder",                     "1 pack",    40.4,               0.6,                  1.2,                  6.8),
    @(306, "tomato millet with powder",                  "1 serving", 323.3333333333333,  7.5,                  11.66666666666667,    42.33333333333334),
    @(307, "coffee with oat creamer",                    "1 serving", 50,                 0,                    2,                    8),
    @(308, "red apple pear orange smoothie no yogurt",   "1 serving", 253.75,             1.125,                1.25,                 58.625),
    @(309, "mtr lemon rice powder",                      "1 pack",    57,                 0.8,                  2.6,                  7.6),
    @(310, "capsicum besan sabji",                       "1 serving", 263.48,             19.18,                9.98,                 25.24)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
